$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.815.86"
$ws.Range("E2").Value = "  +0.51%  "
$ws.Range("D3").Value = "2.293.23"
$ws.Range("E3").Value = "  -0.44%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.07"
$ws.Range("E5").Value = "  -0.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.10"
$ws.Range("E6").Value = "  -0.14%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.625"
$ws.Range("E7").Value = "  -1.01%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("E9").Value = "  -1.76%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.43"
$ws.Range("E10").Value = "  -1.79%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0904"
$ws.Range("E11").Value = "  -0.79%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.47"
$ws.Range("E12").Value = "  +1.41%  "
$ws.Range("E13").Value = "  +2.10%  "
$ws.Range("E14").Value = "  +4.25%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.31"
$ws.Range("E15").Value = "  +0.24%  "
$ws.Range("D16").Value = "2.642.54"
$ws.Range("E16").Value = "  -0.28%  "
$ws.Range("D17").Value = "2.308.92"
$ws.Range("E17").Value = "  +0.20%  "
$ws.Range("D18").Value = "42.762.39"
$ws.Range("E18").Value = "  +0.79%  "
$ws.Range("E19").Value = "  -0.98%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.92"
$ws.Range("E20").Value = "  +24.95%  "
$ws.Range("E21").Value = "  -0.58%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "74.00"
$ws.Range("E22").Value = "  +1.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.56"
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "263.60"
$ws.Range("E24").Value = "  -4.70%  "
$ws.Range("E25").Value = "  -3.51%  "
$ws.Range("E26").Value = "  +0.40%  "
$ws.Range("E27").Value = "  +0.56%  "
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.36"
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("B29").Value = "Filecoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.07"
$ws.Range("E29").Value = "  +19.95%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "37.50"
$ws.Range("E31").Value = "  +4.62%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "166.51"
$ws.Range("E32").Value = "  +0.98%  "
$ws.Range("E33").Value = "  -0.60%  "
$ws.Range("E34").Value = "  -4.58%  "
$ws.Range("E35").Value = "  -0.57%  "
$ws.Range("E36").Value = "  -1.35%  "
$ws.Range("E37").Value = "  -0.61%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0350"
$ws.Range("E38").Value = "  -6.23%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.81"
$ws.Range("E39").Value = "  +1.63%  "
$ws.Range("E40").Value = "  -2.73%  "
$ws.Range("E41").Value = "  +4.19%  "
$ws.Range("E42").Value = "  +1.13%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "69.47"
$ws.Range("E43").Value = "  -0.94%  "
$ws.Range("E44").Value = "  +0.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.08"
$ws.Range("E45").Value = "  -3.27%  "
$ws.Range("E46").Value = "  +0.51%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "114.10"
$ws.Range("E47").Value = "  +0.95%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "80.35"
$ws.Range("E48").Value = "  -2.81%  "
$ws.Range("D49").Value = "1.721.65"
$ws.Range("E49").Value = "  +8.25%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.76"
$ws.Range("E50").Value = "  -1.44%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.14"
$ws.Range("E51").Value = "  +2.02%  "
